# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E60) is re-populated with the period
# codes in ascending chronological order (1607 .. 2003) instead of the
# previous descending order (2003 .. 1607), and the corresponding
# "Valor Mora" column (F16:F60) values are re-matched to their periods:
# periods 1607-1808 (rows 16-41) now carry 24640, and periods 1809-2003
# (rows 42-60) carry 31249.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ascending list of period codes (yymm) that belong in E16:E60, top to bottom.
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

# New "Valor Mora" amounts for F16:F60 (same row order as $periods above):
# rows 16-41 (periods 1607-1808) = 24640, rows 42-60 (periods 1809-2003) = 31249.
$valores = @(
    24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
